$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: the author's window was resized (bookViews windowWidth/
# windowHeight 23003x11351 -> 23651x11748) while editing. Not all hosts
# persist window chrome geometry, so these are wrapped defensively.
try { $excel.ActiveWindow.Width = 23651 } catch {}
try { $excel.ActiveWindow.Height = 11748 } catch {}

# Make the new rows (7-12) use the same cell style ("常规"/Normal w/ 微软雅黑 font)
# already used by the existing example rows (4-6), across every column from
# A to J (mirrors the author copying the row format down before filling values).
$fmtRange = $ws.Range("A7:J12")
$fmtRange.Font.Name = "微软雅黑"

# The old sheet carried a stale outlineLevelRow="5" (no row actually used an
# outline level). Touching each row's OutlineLevel clears that stale metadata.
for ($r = 7; $r -le 12; $r++) {
  $ws.Rows.Item($r).OutlineLevel = 0
}

# Row 7: new "list" example table (Examples.TbExampleList)
$ws.Range("B7").Value = "Examples.TbExampleList"
$ws.Range("C7").Value = "ExampleList"
$ws.Range("D7").Value = $true
$ws.Range("E7").Value = "../列表表@_示例.xlsx"
$ws.Range("G7").Value = "list"

# Row 8: section separator, same convention as the existing "##" rows
$ws.Range("A8").Value = "##"

# Row 9: Game.TbCharacterConfig
$ws.Range("B9").Value = "Game.TbCharacterConfig"
$ws.Range("C9").Value = "CharacterConfig"
$ws.Range("D9").Value = $true
$ws.Range("E9").Value = "../J-角色.xlsx"

# Row 10: Game.TbPlayerInitialConfig
$ws.Range("B10").Value = "Game.TbPlayerInitialConfig"
$ws.Range("C10").Value = "PlayerInitialConfig"
$ws.Range("D10").Value = $true
$ws.Range("E10").Value = "../W-玩家初始化.xlsx"

# Row 11: Game.TbMainLevelConfig
$ws.Range("B11").Value = "Game.TbMainLevelConfig"
$ws.Range("C11").Value = "MainLevelConfig"
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = "../Z-主线关卡.xlsx"

# Row 12: Game.TbGlobalConfig (mode = "one", like ExampleSingleton above)
$ws.Range("B12").Value = "Game.TbGlobalConfig"
$ws.Range("C12").Value = "GlobalConfig"
$ws.Range("D12").Value = $true
$ws.Range("E12").Value = "../T-通用全局.xlsx"
$ws.Range("G12").Value = "one"

# Restore the selection to where the author left off editing
$ws.Range("E15").Select()
